$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (column D) values are written as text, matching the source
# data (values like "66.515.39" are not valid numbers, and values like
# "6.68" must stay literal strings rather than being parsed as numbers).
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.515.39'
$ws.Range("E2").Value = '  -4.04%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.340.39'
$ws.Range("E3").Value = '  -0.50%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.54'
$ws.Range("E5").Value = '  -3.05%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.29'
$ws.Range("E6").Value = '  -5.36%  '

# Row 7
$ws.Range("E7").Value = '  +2.97%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.128'
$ws.Range("E9").Value = '  -3.60%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.68'
$ws.Range("E10").Value = '  -1.43%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.403'
$ws.Range("E11").Value = '  -3.64%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.929.51'
$ws.Range("E12").Value = '  -0.38%  '

# Row 13
$ws.Range("E13").Value = '  -1.03%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = '66.654.43'
$ws.Range("E14").Value = '  -3.93%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '26.80'
$ws.Range("E15").Value = '  -6.00%  '

# Row 16
$ws.Range("E16").Value = '  -2.37%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.348.79'
$ws.Range("E17").Value = '  -1.36%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '432.98'
$ws.Range("E18").Value = '  -4.44%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.58'
$ws.Range("E19").Value = '  -1.28%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.67'
$ws.Range("E20").Value = '  -2.68%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.58'
$ws.Range("E21").Value = '  -3.22%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.42'
$ws.Range("E22").Value = '  -3.54%  '

# Row 23
$ws.Range("E23").Value = '  +0.12%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.518'
$ws.Range("E24").Value = '  -1.14%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000116'
$ws.Range("E25").Value = '  -4.14%  '

# Row 26
$ws.Range("E26").Value = '  +1.20%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.04'
$ws.Range("E27").Value = '  -4.26%  '

# Row 28
$ws.Range("E28").Value = '  -0.15%  '

# Row 29
$ws.Range("E29").Value = '  -2.79%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.82'
$ws.Range("E30").Value = '  -2.19%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("B31").Value = 'USDe'
$ws.Range("C31").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.04%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = '5.25'
$ws.Range("E32").Value = '  -5.74%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.77'
$ws.Range("E33").Value = '  -3.02%  '

# Row 34
$ws.Range("E34").Value = '  -5.32%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.94'
$ws.Range("E35").Value = '  -2.54%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.47'
$ws.Range("E36").Value = '  -5.94%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.87'
$ws.Range("E37").Value = '  +1.97%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.79'
$ws.Range("E38").Value = '  -7.62%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.810.84'
$ws.Range("E39").Value = '  +3.66%  '

# Row 40
$ws.Range("E40").Value = '  -1.06%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.43'
$ws.Range("E41").Value = '  -3.87%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.21'
$ws.Range("E42").Value = '  -4.64%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.30'
$ws.Range("E43").Value = '  -1.61%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0667'
$ws.Range("E44").Value = '  -3.33%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.31'
$ws.Range("E45").Value = '  -4.43%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.33'
$ws.Range("E46").Value = '  -6.87%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '325.13'
$ws.Range("E47").Value = '  -3.25%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0272'
$ws.Range("E48").Value = '  -4.37%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.102'
$ws.Range("E49").Value = '  +0.84%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.973'
$ws.Range("E50").Value = '  -3.80%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.15'
$ws.Range("E51").Value = '  -2.78%  '
